$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$story = "He sits alone in a place desolate and unwanted. He weeps over his beloved wife. He weeps over his beloved son. Nothing can silence the pain and the sorrow in his now frozen heart. He rules over the ruins of the old Christmas Tree Lot. He calls out to the creatures of ice. He seeks her, he seeks him. His rage cannot be quenched."

$ws.Range("A5").Value = "The Frozen King"
$ws.Range("B5").Value = $story
$ws.Range("C5").Value = "frozen-king"
$ws.Range("D5").Value = "The Frozen King"
$ws.Range("E5").Value = "Frozen Child of Fear,Corrupted Christmas Tree,Wailing Banshee of Ice,Bloody Snowman of rage,Faithless Prince of the Snow Garden,Zombified Cat of Yesterday,Faithless Priest of The Old Church"
$ws.Range("F5").Value = "Forgotten Christmas Tree Lot"
$ws.Range("G5").Value = "Dilapidated House of the Drunk,Forgotten Christmas Tree Lot,Banshee Fields of Tomorrow,Ice Pirates Hideout"
$ws.Range("H5").Value = "Corrupted Ice"
$ws.Range("I5").Value = "Ancestral Soldiers Statue"

$ws.Columns.Item(5).ColumnWidth = 224.088
$ws.Columns.Item(6).ColumnWidth = 34.135
$ws.Columns.Item(7).ColumnWidth = 126.112
